# The sheet's column E ("Direccion de red IP" / real-looking values column)
# is selected in full and then hidden, matching the author's edit of
# filling in the DNS table with placeholder ("no reales") values and
# tucking the column away from view.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Columns("E")
[void]$col.Select()
$col.Hidden = $true
